$d = $word.ActiveDocument

# Replace "a basic computer" with "an office computer" in the relevant paragraph
$d.Content.Find.Execute("hundred for a basic computer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "hundred for an office computer", 2)
